$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1234
$ws.Range("C3").Value = "TESTENOVO"
$ws.Range("D3").Value = "CASSA"
$ws.Range("E3").Value = "SAUDE"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "teste20205"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "456654"
$ws.Range("D4").Value = "pc"
$ws.Range("E4").Value = "cpd"

$ws.Range("A5").Value = 4
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "78963"
$ws.Range("C5").Value = "ULTIMO"
$ws.Range("D5").Value = "TESTE"
$ws.Range("E5").Value = "PATRIMONIO"
